# The workbook has a single sheet ("Sheet1") with an A1:B4 table of
# strings. Rows 2-4 previously held a mix of placeholder/credential-looking
# strings; the commit blanks them out to a single space character, leaving
# only the header row (A1/B1) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:B4").Value = " "
